$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 0: the old _GoBack bookmark (after "Two quick & dirty examples:")
# moves up to the new paragraph inserted in Edit 1 below, so remove the
# original one first (before a replacement is created) to avoid any
# ambiguity between the two same-named bookmarks.
# ---------------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# ---------------------------------------------------------------------------
# Edit 1: paragraph ending "...(open source C# compiler)" -- drop the
# trailing manual line break and add a new "Josh looking..." paragraph
# after it (with a relocated _GoBack bookmark at its end).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(5)
$breakPos1 = $p1.Range.End - 2
$d.Range($breakPos1, $breakPos1 + 1).Delete()

$p1 = $d.Paragraphs.Item(5)
$p1.Range.InsertParagraphAfter() | Out-Null
$newPara1 = $d.Paragraphs.Item(6)
$ir1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:b/><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve">Josh </w:t></w:r><w:r><w:rPr><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t>looking @ this, may use Mono or Roslyn.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$ir1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: paragraph ending "...(use a list box that calls/stores  file
# paths)" -- drop the trailing manual line break and add a new "Josh"
# paragraph after it.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(7)
$breakPos2 = $p2.Range.End - 2
$d.Range($breakPos2, $breakPos2 + 1).Delete()

$p2 = $d.Paragraphs.Item(7)
$p2.Range.InsertParagraphAfter() | Out-Null
$newPara2 = $d.Paragraphs.Item(8)
$ir2 = $d.Range($newPara2.Range.Start, $newPara2.Range.End)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t>Josh</w:t></w:r><w:r><w:rPr><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> looking @ this,</w:t></w:r><w:r><w:rPr><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> may end up using tree hierarchy instead</w:t></w:r><w:r><w:rPr><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$ir2.InsertXML($xml2)

Write-Output "done"
